# AO-Smith-MSFT-MCHP-r6.pptx edit
# "Refractor FRAME to DTI"
#
# 1) Bump the cached date-placeholder text on the Handout Master & Notes
#    Master from 11/2/21 -> 11/3/21 (datetimeFigureOut fields).
# 2) Slide 1: update the sub-title/author placeholders -
#      "October 2021"  -> "Wireless Specialist (PACNW)" (+ a trailing blank line)
#      "JR & RW"        -> "Randy Wu, Principal ESE"

$p = $ppt.ActivePresentation

# --- 1. Date placeholders (Handout Master + Notes Master) ---------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.UseFormat = 0
$hm.HeadersFooters.DateAndTime.Value = "11/3/21"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.UseFormat = 0
$nm.HeadersFooters.DateAndTime.Value = "11/3/21"

# --- 2. Slide 1 text updates ---------------------------------------------
$s1 = $p.Slides.Item(1)

# "Text Placeholder 2" (ph idx="12"): October 2021 -> role/title line,
# keeping the original trailing empty paragraph.
$rolePlaceholder = $s1.Shapes.Item(2)
$rolePlaceholder.TextFrame.TextRange.Text = "Wireless Specialist (PACNW)" + [char]13

# "Text Placeholder 3" (ph idx="13"): JR & RW -> full name/title.
$namePlaceholder = $s1.Shapes.Item(3)
$namePlaceholder.TextFrame.TextRange.Text = "Randy Wu, Principal ESE"
